$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the data rows (2 through 8) with the literal column-key placeholders,
# matching the splicing of the database field names into the order rows.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = "sku"
    $ws.Cells.Item($r, 2).Value = "name"
    $ws.Cells.Item($r, 3).Value = "quantity"
    $ws.Cells.Item($r, 4).Value = "cost_per"
    $ws.Cells.Item($r, 5).Value = "total_cost"
}
